# Generate Report for Handoff
# For the files that are "Ready for handoff" (rows 4-7 on each locale sheet),
# bump the Priority to "ht" and refresh the Latest Handoff Datetime.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-09-05 22:40:57"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-09-05 22:41:05"

    # The Overview sheet's "Latest HO Xliff Generate Date" column shares the
    # same underlying text as de-de's handoff datetime, so it refreshes too.
    $overview.Range("G$r").Value = "2016-09-05 22:41:05"
}
